$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.006403796807916251
$ws.Range("D2").Value = 0.230110852657134
$ws.Range("E2").Value = 0.1733927365998476
$ws.Range("F2").Value = 1.119256496137453
$ws.Range("G2").Value = 0.6172754033872252
$ws.Range("H2").Value = 0.6451706849114487
$ws.Range("I2").Value = 0.5624574055654961
$ws.Range("J2").Value = 0.1800749562412989
$ws.Range("N2").Value = 2.873581687696344
$ws.Range("O2").Value = 2.506829514151661
$ws.Range("C3").Value = 0.00560789120331151
$ws.Range("D3").Value = 0.2288708407612177
$ws.Range("E3").Value = 0.1704818582526286
$ws.Range("F3").Value = 1.08621936351328
$ws.Range("G3").Value = 0.5855465827598323
$ws.Range("H3").Value = 0.6344879103896659
$ws.Range("I3").Value = 0.5403325584105687
$ws.Range("J3").Value = 0.1749525151808271
$ws.Range("N3").Value = 2.562605684679454
$ws.Range("O3").Value = 2.416823902178777
$ws.Range("C4").Value = 0.005117407016864917
$ws.Range("D4").Value = 0.2282076702264249
$ws.Range("E4").Value = 0.1687872790564349
$ws.Range("F4").Value = 1.06658017747624
$ws.Range("G4").Value = 0.5664367015796472
$ws.Range("H4").Value = 0.6282672455485283
$ws.Range("I4").Value = 0.5270823313488435
$ws.Range("J4").Value = 0.1719137970435156
$ws.Range("N4").Value = 2.371325805375761
$ws.Range("O4").Value = 2.363042203959907
$ws.Range("C5").Value = 0.004917077860181251
$ws.Range("D5").Value = 0.2279621468700697
$ws.Range("E5").Value = 0.1681200324150822
$ws.Range("F5").Value = 1.058739344460022
$ws.Range("G5").Value = 0.5587425543045583
$ws.Range("H5").Value = 0.6258174501460019
$ws.Range("I5").Value = 0.5217667888827435
$ws.Range("J5").Value = 0.1707022290644034
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 2.341498059993626
$ws.Range("C6").Value = 0.004883785960444698
$ws.Range("D6").Value = 0.2279228720741955
$ws.Range("E6").Value = 0.1680106442375404
$ws.Range("F6").Value = 1.057447179963972
$ws.Range("G6").Value = 0.5574705755689706
$ws.Range("H6").Value = 0.6254158088869843
$ws.Range("I6").Value = 0.5208892200426121
$ws.Range("J6").Value = 0.1705026626185386
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 2.337943149268625
$ws.Range("C7").Value = 0.005114707143540898
$ws.Range("D7").Value = 0.2282042588650839
$ws.Range("E7").Value = 0.1687781859671063
$ws.Range("F7").Value = 1.066473776248586
$ws.Range("G7").Value = 0.5663325581709415
$ws.Range("H7").Value = 0.6282338618372876
$ws.Range("I7").Value = 0.5270103039531904
$ws.Range("J7").Value = 0.1718973492214459
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 2.362750145009358
$ws.Range("C8").Value = 0.006129739967860814
$ws.Range("D8").Value = 0.2296629320517383
$ws.Range("E8").Value = 0.1723698189764029
$ws.Range("F8").Value = 1.107731204184702
$ws.Range("G8").Value = 0.606258023174135
$ws.Range("H8").Value = 0.6414170012355669
$ws.Range("I8").Value = 0.5547592695330081
$ws.Range("J8").Value = 0.1782866114387289
$ws.Range("N8").Value = 2.766433886209654
$ws.Range("O8").Value = 2.475487718664709
$ws.Range("C9").Value = 0.008106133743240207
$ws.Range("D9").Value = 0.2333016362265568
$ws.Range("E9").Value = 0.1801494517381812
$ws.Range("F9").Value = 1.19377010631348
$ws.Range("G9").Value = 0.6875148270819125
$ws.Range("H9").Value = 0.6699565692533156
$ws.Range("I9").Value = 0.6118365369392933
$ws.Range("J9").Value = 0.1916635588980853
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 2.708356419458028
$ws.Range("C10").Value = 0.009549956766697676
$ws.Range("D10").Value = 0.2364487660971122
$ws.Range("E10").Value = 0.1863161131496511
$ws.Range("F10").Value = 1.260133648751889
$ws.Range("G10").Value = 0.7490460878954934
$ws.Range("H10").Value = 0.6925673311482115
$ws.Range("I10").Value = 0.6554091411255314
$ws.Range("J10").Value = 0.2020134014543089
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 2.886696899377966
$ws.Range("C11").Value = 0.01020508346777405
$ws.Range("D11").Value = 0.23798325338295
$ws.Range("E11").Value = 0.1892198975435875
$ws.Range("F11").Value = 1.291013430937312
$ws.Range("G11").Value = 0.7774419290807373
$ws.Range("H11").Value = 0.7032114214596277
$ws.Range("I11").Value = 0.6755907234695968
$ws.Range("J11").Value = 0.2068362859218951
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 2.969418359925953
$ws.Range("C12").Value = 0.01045292532636211
$ws.Range("D12").Value = 0.2385790925780924
$ws.Range("E12").Value = 0.1903336766724024
$ws.Range("F12").Value = 1.302806358558584
$ws.Range("G12").Value = 0.7882532541530907
$ws.Range("H12").Value = 0.7072936193178236
$ws.Range("I12").Value = 0.6832849390733458
$ws.Range("J12").Value = 0.2086791464770243
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 3.000972727489909
$ws.Range("C13").Value = 0.01039955883162946
$ws.Range("D13").Value = 0.2384501118139468
$ws.Range("E13").Value = 0.1900931736462894
$ws.Range("F13").Value = 1.300262116966209
$ws.Range("G13").Value = 0.7859222396498353
$ws.Range("H13").Value = 0.7064121542274222
$ws.Range("I13").Value = 0.6816255428246905
$ws.Range("J13").Value = 0.2082815171670092
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 2.994166714265816
$ws.Range("C14").Value = 0.01022547838171306
$ws.Range("D14").Value = 0.238031977726834
$ws.Range("E14").Value = 0.1893112446014342
$ws.Range("F14").Value = 1.291981648967194
$ws.Range("G14").Value = 0.7783302109860131
$ws.Range("H14").Value = 0.7035462340550964
$ws.Range("I14").Value = 0.6762226907580953
$ws.Range("J14").Value = 0.2069875674365846
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 2.97200975396953
$ws.Range("C15").Value = 0.01011881778805446
$ws.Range("D15").Value = 0.2377777802140173
$ws.Range("E15").Value = 0.188834137283223
$ws.Range("F15").Value = 1.286922574371943
$ws.Range("G15").Value = 0.7736874903795297
$ws.Range("H15").Value = 0.7017974830692708
$ws.Range("I15").Value = 0.6729200480567812
$ws.Range("J15").Value = 0.2061971411045818
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 2.958467881712068
$ws.Range("C16").Value = 0.00950710902657903
$ws.Range("D16").Value = 0.2363505509288757
$ws.Range("E16").Value = 0.1861283276523835
$ws.Range("F16").Value = 1.258129489543307
$ws.Range("G16").Value = 0.7471985191326098
$ws.Range("H16").Value = 0.6918789240660317
$ws.Range("I16").Value = 0.6540974831536914
$ws.Range("J16").Value = 0.2017005264407175
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 2.881322956879842
$ws.Range("C17").Value = 0.009131416756979149
$ws.Range("D17").Value = 0.2355013134694843
$ws.Range("E17").Value = 0.1844936458195647
$ws.Range("F17").Value = 1.240642838012093
$ws.Range("G17").Value = 0.7310522582826877
$ws.Range("H17").Value = 0.6858859766674641
$ws.Range("I17").Value = 0.6426427683739604
$ws.Range("J17").Value = 0.198971406821812
$ws.Range("N17").Value = 3.94211849064385
$ws.Range("O17").Value = 2.834405460557264
$ws.Range("C18").Value = 0.008915170754271173
$ws.Range("D18").Value = 0.2350225360989384
$ws.Range("E18").Value = 0.1835626944836335
$ws.Range("F18").Value = 1.230649981622676
$ws.Range("G18").Value = 0.7218034608712856
$ws.Range("H18").Value = 0.6824727215651478
$ws.Range("I18").Value = 0.6360882161767165
$ws.Range("J18").Value = 0.197412481688545
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 2.807569719244327
$ws.Range("C19").Value = 0.008841926462054062
$ws.Range("D19").Value = 0.2348620939325912
$ws.Range("E19").Value = 0.1832490829698941
$ws.Range("F19").Value = 1.227277733275628
$ws.Range("G19").Value = 0.7186785153292305
$ws.Range("H19").Value = 0.6813228456936713
$ws.Range("I19").Value = 0.6338747763726786
$ws.Range("J19").Value = 0.1968865085890883
$ws.Range("N19").Value = 3.828614786363971
$ws.Range("O19").Value = 2.798509342108559
$ws.Range("C20").Value = 0.009171426163170793
$ws.Range("D20").Value = 0.2355907144344798
$ws.Range("E20").Value = 0.1846667005584735
$ws.Range("F20").Value = 1.242497593576758
$ws.Range("G20").Value = 0.7327671101403723
$ws.Range("H20").Value = 0.6865204453715421
$ws.Range("I20").Value = 0.6438586334458734
$ws.Range("J20").Value = 0.1992608089660024
$ws.Range("N20").Value = 3.957806003281064
$ws.Range("O20").Value = 2.839384384731829
$ws.Range("C21").Value = 0.01027661655066225
$ws.Range("D21").Value = 0.2381543934126569
$ws.Range("E21").Value = 0.1895405312281895
$ws.Range("F21").Value = 1.294411123941259
$ws.Range("G21").Value = 0.7805585862729743
$ws.Range("H21").Value = 0.7043866263950349
$ws.Range("I21").Value = 0.6778082303904114
$ws.Range("J21").Value = 0.2073671826047132
$ws.Range("N21").Value = 4.391158149571311
$ws.Range("O21").Value = 2.978511557254706
$ws.Range("C22").Value = 0.01099752128996556
$ws.Range("D22").Value = 0.2399159389027687
$ws.Range("E22").Value = 0.192808511755473
$ws.Range("F22").Value = 1.32891918798984
$ws.Range("G22").Value = 0.812133784576929
$ws.Range("H22").Value = 0.7163634517991682
$ws.Range("I22").Value = 0.7002987435781307
$ws.Range("J22").Value = 0.2127615894427208
$ws.Range("N22").Value = 4.673791817957863
$ws.Range("O22").Value = 3.070777365247523
$ws.Range("C23").Value = 0.01061288889831502
$ws.Range("D23").Value = 0.2389679055612106
$ws.Range("E23").Value = 0.1910567634459497
$ws.Range("F23").Value = 1.310448518602954
$ws.Range("G23").Value = 0.7952502651589839
$ws.Range("H23").Value = 0.709943727824367
$ws.Range("I23").Value = 0.6882674244892684
$ws.Range("J23").Value = 0.2098736547900302
$ws.Range("N23").Value = 4.523002190005457
$ws.Range("O23").Value = 3.021410791055587
$ws.Range("C24").Value = 0.009153338706482828
$ws.Range("D24").Value = 0.2355502667882376
$ws.Range("E24").Value = 0.184588434949049
$ws.Range("F24").Value = 1.241658870283544
$ws.Range("G24").Value = 0.7319917200247232
$ws.Range("H24").Value = 0.6862335018541046
$ws.Range("I24").Value = 0.6433088445505888
$ws.Range("J24").Value = 0.1991299388575101
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 2.837132984039215
$ws.Range("C25").Value = 0.007572937895453435
$ws.Range("D25").Value = 0.2322340231463187
$ws.Range("E25").Value = 0.1779658061177756
$ws.Range("F25").Value = 1.169942548261744
$ws.Range("G25").Value = 0.6652128538278532
$ws.Range("H25").Value = 0.6619476978387411
$ws.Range("I25").Value = 0.5961090942241327
$ws.Range("J25").Value = 0.1879534732702695
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 2.644090394289606
